$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4320.14475225888
$ws.Range("C3").Value = 4320.14475225888
$ws.Range("C4").Value = 4262.024763711572
$ws.Range("C5").Value = 4113.425971195736
$ws.Range("C6").Value = 4113.425971195736
$ws.Range("C7").Value = 4113.425971195736
$ws.Range("C8").Value = 3995.127283963309
$ws.Range("C9").Value = 3976.960410242727
$ws.Range("C10").Value = 3976.960410242727
$ws.Range("C11").Value = 3976.960410242727
$ws.Range("C12").Value = 3976.960410242727
